$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-02-07 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-02-08 Thursday", 2) | Out-Null
$d.Content.Find.Execute("536÷4=134, 0", $true, $false, $false, $false, $false, $true, 1, $false, "799÷3=266, 1", 2) | Out-Null
$d.Content.Find.Execute("510÷7=72, 6", $true, $false, $false, $false, $false, $true, 1, $false, "579÷8=72, 3", 2) | Out-Null
$d.Content.Find.Execute("237÷3=79, 0", $true, $false, $false, $false, $false, $true, 1, $false, "870÷8=108, 6", 2) | Out-Null
$d.Content.Find.Execute("112÷3=37, 1", $true, $false, $false, $false, $false, $true, 1, $false, "455÷6=75, 5", 2) | Out-Null
$d.Content.Find.Execute("630÷6=105, 0", $true, $false, $false, $false, $false, $true, 1, $false, "919÷7=131, 2", 2) | Out-Null
$d.Content.Find.Execute("968÷9=107, 5", $true, $false, $false, $false, $false, $true, 1, $false, "425÷2=212, 1", 2) | Out-Null
$d.Content.Find.Execute("301÷3=100, 1", $true, $false, $false, $false, $false, $true, 1, $false, "294÷8=36, 6", 2) | Out-Null
$d.Content.Find.Execute("263÷7=37, 4", $true, $false, $false, $false, $false, $true, 1, $false, "741÷6=123, 3", 2) | Out-Null
$d.Content.Find.Execute("965÷8=120, 5", $true, $false, $false, $false, $false, $true, 1, $false, "448÷7=64, 0", 2) | Out-Null
$d.Content.Find.Execute("481÷7=68, 5", $true, $false, $false, $false, $false, $true, 1, $false, "284÷8=35, 4", 2) | Out-Null
$d.Content.Find.Execute("690÷3=230, 0", $true, $false, $false, $false, $false, $true, 1, $false, "577÷6=96, 1", 2) | Out-Null
$d.Content.Find.Execute("881÷5=176, 1", $true, $false, $false, $false, $false, $true, 1, $false, "135÷4=33, 3", 2) | Out-Null
$d.Content.Find.Execute("471÷7=67, 2", $true, $false, $false, $false, $false, $true, 1, $false, "397÷3=132, 1", 2) | Out-Null
$d.Content.Find.Execute("379÷5=75, 4", $true, $false, $false, $false, $false, $true, 1, $false, "467÷7=66, 5", 2) | Out-Null
$d.Content.Find.Execute("299÷8=37, 3", $true, $false, $false, $false, $false, $true, 1, $false, "647÷7=92, 3", 2) | Out-Null
$d.Content.Find.Execute("698÷5=139, 3", $true, $false, $false, $false, $false, $true, 1, $false, "375÷9=41, 6", 2) | Out-Null
$d.Content.Find.Execute("819÷2=409, 1", $true, $false, $false, $false, $false, $true, 1, $false, "914÷4=228, 2", 2) | Out-Null
$d.Content.Find.Execute("378÷9=42, 0", $true, $false, $false, $false, $false, $true, 1, $false, "585÷2=292, 1", 2) | Out-Null
$d.Content.Find.Execute("179÷5=35, 4", $true, $false, $false, $false, $false, $true, 1, $false, "133÷5=26, 3", 2) | Out-Null
$d.Content.Find.Execute("900÷8=112, 4", $true, $false, $false, $false, $false, $true, 1, $false, "207÷9=23, 0", 2) | Out-Null
$d.Content.Find.Execute("400÷7=57, 1", $true, $false, $false, $false, $false, $true, 1, $false, "682÷7=97, 3", 2) | Out-Null
$d.Content.Find.Execute("180÷4=45, 0", $true, $false, $false, $false, $false, $true, 1, $false, "321÷5=64, 1", 2) | Out-Null
$d.Content.Find.Execute("652÷3=217, 1", $true, $false, $false, $false, $false, $true, 1, $false, "512÷2=256, 0", 2) | Out-Null
$d.Content.Find.Execute("851÷6=141, 5", $true, $false, $false, $false, $false, $true, 1, $false, "703÷4=175, 3", 2) | Out-Null
$d.Content.Find.Execute("484÷9=53, 7", $true, $false, $false, $false, $false, $true, 1, $false, "882÷6=147, 0", 2) | Out-Null
